{"js": "// Body: \"A QWREW,\" -> \"A QWR,\"  (bold placeholder run in the body text)\nconst bodySearch = context.document.body.search(\"QWREW\", {\n  matchCase: true,\n  matchWholeWord: true\n});\nbodySearch.load(\"items\");\nawait context.sync();\nbodySearch.items.forEach((r) => r.insertText(\"QWR\", Word.InsertLocation.replace));\nawait context.sync();\n\n// Header: grab the primary header of every section (the doc only has one,\n// but loop defensively in case a template ever grows more sections).\nconst sections = context.document.sections;\nsections.load(\"items\");\nawait context.sync();\n\nfor (const sec of sections.items) {\n  const header = sec.getHeader(\"Primary\");\n\n  // \"DIRETORIA DE ENSINO REGIAO REW\" -> \"... QWER\"\n  const regiaoSearch = header.search(\"REW\", { matchCase: true, matchWholeWord: true });\n  regiaoSearch.load(\"items\");\n\n  // \"QWREW \u2013 DEP.\" -> \"QWR \u2013 DEP.\"\n  const depSearch = header.search(\"QWREW\", { matchCase: true, matchWholeWord: true });\n  depSearch.load(\"items\");\n\n  // Address line: five \"Rew\" runs -> \"Qwer\"\n  const rewSearch = header.search(\"Rew\", { matchCase: true, matchWholeWord: true });\n  rewSearch.load(\"items\");\n\n  // CEP / Tel / Email lines: three \"rew\" runs -> \"qwer\"\n  const lowerRewSearch = header.search(\"rew\", { matchCase: true, matchWholeWord: true });\n  lowerRewSearch.load(\"items\");\n\n  await context.sync();\n\n  regiaoSearch.items.forEach((r) => r.insertText(\"QWER\", Word.InsertLocation.replace));\n  depSearch.items.forEach((r) => r.insertText(\"QWR\", Word.InsertLocation.replace));\n  rewSearch.items.forEach((r) => r.insertText(\"Qwer\", Word.InsertLocation.replace));\n  lowerRewSearch.items.forEach((r) => r.insertText(\"qwer\", Word.InsertLocation.replace));\n\n  await context.sync();\n}\n", "ps1": "# RBA 2.3 - Relatorio e Email\n# Tidy up the placeholder \"QWREW\"/\"REW\"/\"Rew\"/\"rew\" tokens that were left in\n# the convocation template (body salutation + the address block in the page\n# header).\n\n$d = $word.ActiveDocument\n\nfunction Replace-AllText($range, [string]$findText, [string]$replaceText) {\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.Text = $replaceText\n    $find.Forward = $true\n    $find.Wrap = 0          # wdFindStop - do not wrap past the end of range\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $true\n    $find.MatchWildcards = $false\n    $find.MatchSoundsLike = $false\n    $find.MatchAllWordForms = $false\n    $find.Execute([ref]$findText, [ref]$true, [ref]$true, [ref]$false, [ref]$false, `\n        [ref]$false, [ref]$true, [ref]0, [ref]$false, [ref]$replaceText, [ref]2) | Out-Null\n}\n\n# Body salutation: \"A QWREW,\" -> \"A QWR,\"\nReplace-AllText $d.Content \"QWREW\" \"QWR\"\n\n# Page header (single section in this template).\nforeach ($section in $d.Sections) {\n    $headerRange = $section.Headers(1).Range   # wdHeaderFooterPrimary\n\n    # \"DIRETORIA DE ENSINO REGIAO REW\" -> \"... QWER\"\n    Replace-AllText $headerRange \"REW\" \"QWER\"\n\n    # \"QWREW - DEP.\" -> \"QWR - DEP.\"\n    Replace-AllText $headerRange \"QWREW\" \"QWR\"\n\n    # Address line: five \"Rew\" runs -> \"Qwer\"\n    Replace-AllText $headerRange \"Rew\" \"Qwer\"\n\n    # CEP / Tel / Email lines: three \"rew\" runs -> \"qwer\"\n    Replace-AllText $headerRange \"rew\" \"qwer\"\n}\n"}
